{"js": "// The document body is: one paragraph (date header) followed by a single\n// table. The table has 20 rows x 5 columns; only every 4th row (0, 4, 8,\n// 12, 16 - zero based) actually holds a division-problem answer such as\n// \"24\u00f75=4, 4\" as plain text in a single run inside a single paragraph in\n// each of its 5 cells. The remaining rows are blank spacer rows.\n//\n// We walk the table row by row / cell by cell and, for the five rows that\n// contain text, replace the text of that cell's single paragraph in place\n// (preserving the existing run/paragraph formatting) with the new value.\n//\n// We key each replacement off its exact row/column position (rather than\n// searching the whole document for the old text and replacing the first\n// match) because some of the new values are identical to *other* cells'\n// old values - a plain text search/replace would be ambiguous/order\n// dependent. We also sanity check the existing text before overwriting it.\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Map of row index -> array of 5 {before, after} pairs (document order,\n// left to right), matching the diff.\nconst rowReplacements = {\n  0: [\n    { before: \"24\u00f75=4, 4\", after: \"50\u00f73=16, 2\" },\n    { before: \"41\u00f72=20, 1\", after: \"99\u00f73=33, 0\" },\n    { before: \"75\u00f79=8, 3\", after: \"89\u00f74=22, 1\" },\n    { before: \"68\u00f78=8, 4\", after: \"42\u00f72=21, 0\" },\n    { before: \"32\u00f78=4, 0\", after: \"65\u00f78=8, 1\" },\n  ],\n  4: [\n    { before: \"21\u00f79=2, 3\", after: \"83\u00f73=27, 2\" },\n    { before: \"33\u00f79=3, 6\", after: \"54\u00f76=9, 0\" },\n    { before: \"77\u00f75=15, 2\", after: \"60\u00f75=12, 0\" },\n    { before: \"45\u00f78=5, 5\", after: \"97\u00f79=10, 7\" },\n    { before: \"72\u00f75=14, 2\", after: \"50\u00f72=25, 0\" },\n  ],\n  8: [\n    { before: \"67\u00f79=7, 4\", after: \"81\u00f78=10, 1\" },\n    { before: \"68\u00f77=9, 5\", after: \"32\u00f72=16, 0\" },\n    { before: \"49\u00f77=7, 0\", after: \"61\u00f75=12, 1\" },\n    { before: \"50\u00f75=10, 0\", after: \"19\u00f78=2, 3\" },\n    { before: \"10\u00f76=1, 4\", after: \"58\u00f77=8, 2\" },\n  ],\n  12: [\n    { before: \"62\u00f74=15, 2\", after: \"33\u00f74=8, 1\" },\n    { before: \"20\u00f73=6, 2\", after: \"19\u00f79=2, 1\" },\n    { before: \"83\u00f73=27, 2\", after: \"56\u00f76=9, 2\" },\n    { before: \"67\u00f77=9, 4\", after: \"46\u00f73=15, 1\" },\n    { before: \"97\u00f73=32, 1\", after: \"25\u00f76=4, 1\" },\n  ],\n  16: [\n    { before: \"62\u00f78=7, 6\", after: \"67\u00f74=16, 3\" },\n    { before: \"69\u00f78=8, 5\", after: \"27\u00f72=13, 1\" },\n    { before: \"56\u00f73=18, 2\", after: \"46\u00f72=23, 0\" },\n    { before: \"85\u00f77=12, 1\", after: \"47\u00f75=9, 2\" },\n    { before: \"82\u00f79=9, 1\", after: \"11\u00f77=1, 4\" },\n  ],\n};\n\n// Load the cells for every row that needs editing.\nconst targetRowIndexes = Object.keys(rowReplacements).map(Number);\nconst cellsByRow = {};\nfor (const r of targetRowIndexes) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  cellsByRow[r] = cells;\n}\nawait context.sync();\n\n// Load each cell's paragraph (text) so we can verify + edit it.\nconst paraByRowCol = {};\nfor (const r of targetRowIndexes) {\n  const cells = cellsByRow[r];\n  for (let c = 0; c < cells.items.length; c++) {\n    const para = cells.items[c].body.paragraphs.getFirst();\n    para.load(\"text\");\n    paraByRowCol[r + \"_\" + c] = para;\n  }\n}\nawait context.sync();\n\n// Verify existing text then queue up the replacement.\nfor (const r of targetRowIndexes) {\n  const pairs = rowReplacements[r];\n  for (let c = 0; c < pairs.length; c++) {\n    const { before, after } = pairs[c];\n    const para = paraByRowCol[r + \"_\" + c];\n    if (para.text !== before) {\n      throw new Error(\n        \"Unexpected text in row \" +\n          r +\n          \" col \" +\n          c +\n          \": \" +\n          JSON.stringify(para.text) +\n          \" (expected \" +\n          JSON.stringify(before) +\n          \")\"\n      );\n    }\n    para.insertText(after, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 division-problem answers in the single table of the\n# document. The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17\n# (1-based) hold answers such as \"24\u00f75=4, 4\" - the remaining rows are\n# blank spacer rows.\n#\n# We address each cell by its exact (row, column) position - rather than\n# a document-wide Find/Replace - because a couple of the new values are\n# identical to *other* cells' old values, which would make a plain\n# search-and-replace pass order dependent / ambiguous. We also verify the\n# existing text before overwriting it.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word.Range.Text on a table cell always ends with a CR + BEL \"cell mark\"\n# pair; strip those two characters off before comparing against our\n# expected plain-text values.\n$cr = [char]13\n$bel = [char]7\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Before = \"24\u00f75=4, 4\"; After = \"50\u00f73=16, 2\" }\n    @{ Row = 1; Col = 2; Before = \"41\u00f72=20, 1\"; After = \"99\u00f73=33, 0\" }\n    @{ Row = 1; Col = 3; Before = \"75\u00f79=8, 3\"; After = \"89\u00f74=22, 1\" }\n    @{ Row = 1; Col = 4; Before = \"68\u00f78=8, 4\"; After = \"42\u00f72=21, 0\" }\n    @{ Row = 1; Col = 5; Before = \"32\u00f78=4, 0\"; After = \"65\u00f78=8, 1\" }\n    @{ Row = 5; Col = 1; Before = \"21\u00f79=2, 3\"; After = \"83\u00f73=27, 2\" }\n    @{ Row = 5; Col = 2; Before = \"33\u00f79=3, 6\"; After = \"54\u00f76=9, 0\" }\n    @{ Row = 5; Col = 3; Before = \"77\u00f75=15, 2\"; After = \"60\u00f75=12, 0\" }\n    @{ Row = 5; Col = 4; Before = \"45\u00f78=5, 5\"; After = \"97\u00f79=10, 7\" }\n    @{ Row = 5; Col = 5; Before = \"72\u00f75=14, 2\"; After = \"50\u00f72=25, 0\" }\n    @{ Row = 9; Col = 1; Before = \"67\u00f79=7, 4\"; After = \"81\u00f78=10, 1\" }\n    @{ Row = 9; Col = 2; Before = \"68\u00f77=9, 5\"; After = \"32\u00f72=16, 0\" }\n    @{ Row = 9; Col = 3; Before = \"49\u00f77=7, 0\"; After = \"61\u00f75=12, 1\" }\n    @{ Row = 9; Col = 4; Before = \"50\u00f75=10, 0\"; After = \"19\u00f78=2, 3\" }\n    @{ Row = 9; Col = 5; Before = \"10\u00f76=1, 4\"; After = \"58\u00f77=8, 2\" }\n    @{ Row = 13; Col = 1; Before = \"62\u00f74=15, 2\"; After = \"33\u00f74=8, 1\" }\n    @{ Row = 13; Col = 2; Before = \"20\u00f73=6, 2\"; After = \"19\u00f79=2, 1\" }\n    @{ Row = 13; Col = 3; Before = \"83\u00f73=27, 2\"; After = \"56\u00f76=9, 2\" }\n    @{ Row = 13; Col = 4; Before = \"67\u00f77=9, 4\"; After = \"46\u00f73=15, 1\" }\n    @{ Row = 13; Col = 5; Before = \"97\u00f73=32, 1\"; After = \"25\u00f76=4, 1\" }\n    @{ Row = 17; Col = 1; Before = \"62\u00f78=7, 6\"; After = \"67\u00f74=16, 3\" }\n    @{ Row = 17; Col = 2; Before = \"69\u00f78=8, 5\"; After = \"27\u00f72=13, 1\" }\n    @{ Row = 17; Col = 3; Before = \"56\u00f73=18, 2\"; After = \"46\u00f72=23, 0\" }\n    @{ Row = 17; Col = 4; Before = \"85\u00f77=12, 1\"; After = \"47\u00f75=9, 2\" }\n    @{ Row = 17; Col = 5; Before = \"82\u00f79=9, 1\"; After = \"11\u00f77=1, 4\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $current = $cell.Range.Text.TrimEnd($cr, $bel)\n    if ($current -ne $item.Before) {\n        throw \"Unexpected text in row $($item.Row) col $($item.Col): [$current] (expected [$($item.Before)])\"\n    }\n    $cell.Range.Text = $item.After\n}\n\n"}
